# Auto-generated Excel COM-interop script applying the scraped diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 38
$ws.Range("H38").Value = 63.5
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()

# ALC row 69
$ws.Range("H69").Value = 7372.048
$ws.Range("J69").Value = 7372.048
$ws.Range("L69").Value = 22116.144
$ws.Range("N69").Value = -23864.144

# ALC row 70
$ws.Range("H70").Value = 2999.25
$ws.Range("I70").Value = 2999
$ws.Range("K70").Value = 8997
$ws.Range("M70").Value = -8727

# ALC row 72
$ws.Range("H72").Value = 7372.048
$ws.Range("J72").Value = 7372.048
$ws.Range("L72").Value = 66348.432
$ws.Range("N72").Value = -75084.432

# ALC row 73
$ws.Range("H73").Value = 2999.25
$ws.Range("I73").Value = 2999
$ws.Range("K73").Value = 8997
$ws.Range("M73").Value = -8061

# ALC row 80
$ws.Range("H80").Value = 514.8125
$ws.Range("I80").Value = 394.1111
$ws.Range("K80").Value = 1182.3333
$ws.Range("M80").Value = -184.3333

# ALC row 83
$ws.Range("H83").Value = 514.8125
$ws.Range("I83").Value = 394.1111
$ws.Range("K83").Value = 3546.9999
$ws.Range("M83").Value = 1445.0001

# ALC row 88
$ws.Range("H88").Value = 2316.2
$ws.Range("J88").Value = 1981.75
$ws.Range("L88").Value = 1981.75
$ws.Range("N88").Value = -2793.75

# ALC row 91
$ws.Range("H91").Value = 2316.2
$ws.Range("J91").Value = 1981.75
$ws.Range("L91").Value = 1981.75
$ws.Range("N91").Value = -4789.75

# ALC row 107
$ws.Range("H107").Value = 404.18182
$ws.Range("I107").Value = 245.5
$ws.Range("J107").Value = 594.6
$ws.Range("K107").Value = 245.5
$ws.Range("L107").Value = 594.6
$ws.Range("M107").Value = 1674.5
$ws.Range("N107").Value = -4434.6

# ALC row 137
$ws.Range("H137").Value = 5449
$ws.Range("I137").Value = 1143
$ws.Range("J137").Value = 13199.8
$ws.Range("K137").Value = 3429
$ws.Range("L137").Value = 39599.39999999999
$ws.Range("M137").Value = -879
$ws.Range("N137").Value = -44699.39999999999

# ALC row 138
$ws.Range("H138").Value = 3202
$ws.Range("I138").Value = 990.25
$ws.Range("J138").Value = 9100
$ws.Range("K138").Value = 2970.75
$ws.Range("L138").Value = 27300
$ws.Range("M138").Value = 2169.25
$ws.Range("N138").Value = -37580

# ALC row 140
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
# ARM row 32
$ws.Range("H32").Value = 3948.48
$ws.Range("I32").Value = 748.1905
$ws.Range("K32").Value = 748.1905
$ws.Range("M32").Value = -461.1905

# ARM row 61
$ws.Range("H61").Value = 1887
$ws.Range("I61").Value = 1349.3334
$ws.Range("K61").Value = 1349.3334
$ws.Range("M61").Value = -1137.3334

# ARM row 74
$ws.Range("H74").Value = 3199.8696
$ws.Range("I74").Value = 2977.3635
$ws.Range("K74").Value = 2977.3635
$ws.Range("M74").Value = -2103.3635

# ARM row 77
$ws.Range("H77").Value = 3199.8696
$ws.Range("I77").Value = 2977.3635
$ws.Range("K77").Value = 14886.8175
$ws.Range("M77").Value = -10518.8175

# ARM row 115
$ws.Range("H115").Value = 20621
$ws.Range("I115").Value = 20621
$ws.Range("J115").Value = 0
$ws.Range("K115").Value = 20621
$ws.Range("L115").Value = 0
$ws.Range("N115").ClearContents()
$ws.Range("M115").Value = -19054

# ARM row 122
$ws.Range("H122").Value = 2125
$ws.Range("I122").Value = 2125
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 6375
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -3925
$ws.Range("N122").ClearContents()

# ARM row 136
$ws.Range("H136").Value = 1887
$ws.Range("I136").Value = 1349.3334
$ws.Range("K136").Value = 4048.0002
$ws.Range("M136").Value = -1498.0002

$ws = $wb.Worksheets.Item("BSM")
# BSM row 43
$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()

# BSM row 99
$ws.Range("H99").Value = 2102.25
$ws.Range("I99").Value = 2010
$ws.Range("K99").Value = 2010
$ws.Range("M99").Value = -512

$ws = $wb.Worksheets.Item("CRP")
# CRP row 18
$ws.Range("H18").Value = 43684
$ws.Range("J18").Value = 43684
$ws.Range("L18").Value = 43684
$ws.Range("N18").Value = -44144

# CRP row 31
$ws.Range("H31").Value = 5364.375
$ws.Range("I31").Value = 3850
$ws.Range("K31").Value = 3850
$ws.Range("M31").Value = -3555

# CRP row 34
$ws.Range("H34").Value = 5364.375
$ws.Range("I34").Value = 3850
$ws.Range("K34").Value = 3850
$ws.Range("M34").Value = -3648

# CRP row 58
$ws.Range("H58").Value = 3256.2
$ws.Range("I58").Value = 1884
$ws.Range("K58").Value = 1884
$ws.Range("M58").Value = -1681

# CRP row 122
$ws.Range("H122").Value = 760.5
$ws.Range("I122").Value = 522
$ws.Range("J122").Value = 999
$ws.Range("K122").Value = 1566
$ws.Range("L122").Value = 2997
$ws.Range("M122").Value = 884
$ws.Range("N122").Value = -7897

# CRP row 134
$ws.Range("H134").Value = 1647.8214
$ws.Range("I134").Value = 922.4167
$ws.Range("K134").Value = 2767.2501
$ws.Range("M134").Value = -232.2501000000002

# CRP row 136
$ws.Range("H136").Value = 3256.2
$ws.Range("I136").Value = 1884
$ws.Range("K136").Value = 5652
$ws.Range("M136").Value = -3102

$ws = $wb.Worksheets.Item("CUL")
# CUL row 92
$ws.Range("H92").Value = 5582.6665
$ws.Range("I92").Value = 2000
$ws.Range("J92").Value = 6299.2
$ws.Range("K92").Value = 6000
$ws.Range("L92").Value = 18897.6
$ws.Range("N92").Value = -21393.6
$ws.Range("M92").Value = -4752

# CUL row 137
$ws.Range("H137").Value = 3457.4546
$ws.Range("I137").Value = 2466.3333
$ws.Range("J137").Value = 3829.125
$ws.Range("K137").Value = 7398.999899999999
$ws.Range("L137").Value = 11487.375
$ws.Range("M137").Value = -2298.999899999999
$ws.Range("N137").Value = -21687.375

$ws = $wb.Worksheets.Item("GSM")
# GSM row 122
$ws.Range("H122").Value = 1772.909
$ws.Range("I122").Value = 1772.909
$ws.Range("K122").Value = 5318.727000000001
$ws.Range("M122").Value = -2868.727000000001

$ws = $wb.Worksheets.Item("LTW")
# LTW row 22
$ws.Range("H22").Value = 1687.4166
$ws.Range("I22").Value = 937.25
$ws.Range("K22").Value = 937.25
$ws.Range("M22").Value = -642.25

# LTW row 27
$ws.Range("H27").Value = 1687.4166
$ws.Range("I27").Value = 937.25
$ws.Range("K27").Value = 937.25
$ws.Range("M27").Value = -830.25

# LTW row 100
$ws.Range("H100").Value = 2314.125
$ws.Range("I100").Value = 2359
$ws.Range("K100").Value = 2359
$ws.Range("M100").Value = -1818

$ws = $wb.Worksheets.Item("WVR")
# WVR row 136
$ws.Range("H136").Value = 1928.44
$ws.Range("I136").Value = 1605.55
$ws.Range("K136").Value = 4816.65
$ws.Range("M136").Value = -2266.65
